$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.02240696193102565
$ws.Range("D2").Value = 0.1747284524113866
$ws.Range("E2").Value = 0.145479424469638
$ws.Range("F2").Value = 1.37011226830522
$ws.Range("G2").Value = 0.002453775259263491
$ws.Range("I2").Value = 0.7068819490681832
$ws.Range("J2").Value = 0.1593150103089762
$ws.Range("K2").Value = 0.9876128801214463
$ws.Range("M2").Value = 0.364592288795599
$ws.Range("N2").Value = 1.349540083739768
$ws.Range("O2").Value = 3.293106829160308
$ws.Range("C3").Value = 0.0202474869338829
$ws.Range("D3").Value = 0.1728109719831679
$ws.Range("E3").Value = 0.1452572531577623
$ws.Range("F3").Value = 1.371808563168543
$ws.Range("G3").Value = 0.002456553525357978
$ws.Range("I3").Value = 0.7066723921334841
$ws.Range("J3").Value = 0.1600267082083562
$ws.Range("K3").Value = 0.884620536743455
$ws.Range("M3").Value = 0.342257071509529
$ws.Range("N3").Value = 1.357518245613896
$ws.Range("O3").Value = 3.303342260291288
$ws.Range("C4").Value = 0.01891199281244127
$ws.Range("D4").Value = 0.1716943500043513
$ws.Range("E4").Value = 0.1451838187859558
$ws.Range("F4").Value = 1.373631225203425
$ws.Range("G4").Value = 0.002458351397402303
$ws.Range("I4").Value = 0.7069904834642031
$ws.Range("J4").Value = 0.1605421419531545
$ws.Range("K4").Value = 0.8213934759581605
$ws.Range("M4").Value = 0.3286431795593145
$ws.Range("N4").Value = 1.362905314122649
$ws.Range("O4").Value = 3.311626771070763
$ws.Range("C5").Value = 0.01836538082999084
$ws.Range("D5").Value = 0.1712546533403199
$ws.Range("E5").Value = 0.1451697543853392
$ws.Range("F5").Value = 1.374570382201377
$ws.Range("G5").Value = 0.002459107249245176
$ws.Range("I5").Value = 0.7072323928676667
$ws.Range("J5").Value = 0.1607719141183779
$ws.Range("K5").Value = 0.7956322273468004
$ws.Range("M5").Value = 0.3231209190688702
$ws.Range("N5").Value = 1.365223588344612
$ws.Range("O5").Value = 3.315505381435258
$ws.Range("C6").Value = 0.018274472689221
$ws.Range("D6").Value = 0.1711825702839533
$ws.Range("E6").Value = 0.1451683776554979
$ws.Range("F6").Value = 1.374738190643697
$ws.Range("G6").Value = 0.002459234161451829
$ws.Range("I6").Value = 0.7072793417990653
$ws.Range("J6").Value = 0.1608112592473248
$ws.Range("K6").Value = 0.7913548984265617
$ws.Range("M6").Value = 0.3222055033299327
$ws.Range("N6").Value = 1.365615970204139
$ws.Range("O6").Value = 3.316179769394665
$ws.Range("C7").Value = 0.01890463065381454
$ws.Range("D7").Value = 0.171688357917148
$ws.Range("E7").Value = 0.1451835648571489
$ws.Range("F7").Value = 1.373643095777219
$ws.Range("G7").Value = 0.002458361497168664
$ws.Range("I7").Value = 0.7069932913793053
$ws.Range("J7").Value = 0.1605451608571222
$ws.Range("K7").Value = 0.8210460312320151
$ws.Range("M7").Value = 0.3285686006526589
$ws.Range("N7").Value = 1.362936080920676
$ws.Range("O7").Value = 3.31167704472773
$ws.Range("C8").Value = 0.02166437130230037
$ws.Range("D8").Value = 0.1740547426080568
$ws.Range("E8").Value = 0.1453897623136555
$ws.Range("F8").Value = 1.370535015200517
$ws.Range("G8").Value = 0.002454714151221416
$ws.Range("I8").Value = 0.7067169323462892
$ws.Range("J8").Value = 0.1595441229283026
$ws.Range("K8").Value = 0.9520998434898615
$ws.Range("M8").Value = 0.3568705545830255
$ws.Range("N8").Value = 1.352189681306513
$ws.Range("O8").Value = 3.29622074231753
$ws.Range("C9").Value = 0.02699971973422777
$ws.Range("D9").Value = 0.1791743952949787
$ws.Range("E9").Value = 0.1462929565641673
$ws.Range("F9").Value = 1.370640340037937
$ws.Range("G9").Value = 0.002448288596374594
$ws.Range("I9").Value = 0.709723470521233
$ws.Range("J9").Value = 0.1582035514335303
$ws.Range("K9").Value = 1.209122278879079
$ws.Range("M9").Value = 0.4131513903177719
$ws.Range("N9").Value = 1.334984303623102
$ws.Range("O9").Value = 3.281795909919254
$ws.Range("C10").Value = 0.03087254093037473
$ws.Range("D10").Value = 0.1832249871114158
$ws.Range("E10").Value = 0.1472597285362198
$ws.Range("F10").Value = 1.37450328087219
$ws.Range("G10").Value = 0.002444006462917782
$ws.Range("I10").Value = 0.714102200030446
$ws.Range("J10").Value = 0.1575982731802554
$ws.Range("K10").Value = 1.397913051383796
$ws.Range("M10").Value = 0.454964031388279
$ws.Range("N10").Value = 1.324692546410816
$ws.Range("O10").Value = 3.280909662011055
$ws.Range("C11").Value = 0.03262409864726123
$ws.Range("D11").Value = 0.1851299125438288
$ws.Range("E11").Value = 0.1477651989221727
$ws.Range("F11").Value = 1.377084028811211
$ws.Range("G11").Value = 0.00244215273727566
$ws.Range("I11").Value = 0.7165668679901103
$ws.Range("J11").Value = 0.1574054013879618
$ws.Range("K11").Value = 1.483777910436743
$ws.Range("M11").Value = 0.4740837870636128
$ws.Range("N11").Value = 1.320518830304991
$ws.Range("O11").Value = 3.282621608996436
$ws.Range("C12").Value = 0.0332858855677074
$ws.Range("D12").Value = 0.1858601550919872
$ws.Range("E12").Value = 0.1479660326716967
$ws.Range("F12").Value = 1.378179775541753
$ws.Range("G12").Value = 0.002441464260038358
$ws.Range("I12").Value = 0.7175682494330928
$ws.Range("J12").Value = 0.1573442272295082
$ws.Range("K12").Value = 1.516288869654659
$ws.Range("M12").Value = 0.4813378695566328
$ws.Range("N12").Value = 1.319011271224014
$ws.Range("O12").Value = 3.283574458496446
$ws.Range("C13").Value = 0.03314342453290919
$ws.Range("D13").Value = 0.1857024900622548
$ws.Range("E13").Value = 0.1479223607825055
$ws.Range("F13").Value = 1.377938516809891
$ws.Range("G13").Value = 0.002441611937081363
$ws.Range("I13").Value = 0.7173495554962344
$ws.Range("J13").Value = 0.1573568745595821
$ws.Range("K13").Value = 1.509287267532102
$ws.Range("M13").Value = 0.4797749634445054
$ws.Range("N13").Value = 1.319332709040779
$ws.Range("O13").Value = 3.283355692651753
$ws.Range("C14").Value = 0.03267857424255283
$ws.Range("D14").Value = 0.1851898122897353
$ws.Range("E14").Value = 0.1477815329075121
$ws.Range("F14").Value = 1.377171801835345
$ws.Range("G14").Value = 0.002442095825786646
$ws.Range("I14").Value = 0.7166478876456068
$ws.Range("J14").Value = 0.1574001308007738
$ws.Range("K14").Value = 1.486452700808343
$ws.Range("M14").Value = 0.4746803095142411
$ws.Range("N14").Value = 1.320393341428229
$ws.Range("O14").Value = 3.282693893998044
$ws.Range("C15").Value = 0.03239364532399236
$ws.Range("D15").Value = 0.1848769378746056
$ws.Range("E15").Value = 0.1476964982861411
$ws.Range("F15").Value = 1.376717597250575
$ws.Range("G15").Value = 0.002442393976460532
$ws.Range("I15").Value = 0.7162269623632724
$ws.Range("J15").Value = 0.1574281713897747
$ws.Range("K15").Value = 1.472465270505722
$ws.Range("M15").Value = 0.4715614785691287
$ws.Range("N15").Value = 1.321052504388959
$ws.Range("O15").Value = 3.282328200047147
$ws.Range("C16").Value = 0.03075786521553425
$ws.Range("D16").Value = 0.1831017433291748
$ws.Range("E16").Value = 0.1472280146662186
$ws.Range("F16").Value = 1.374351198397378
$ws.Range("G16").Value = 0.002444129499409683
$ws.Range("I16").Value = 0.7139506484345617
$ws.Range("J16").Value = 0.1576125375966697
$ws.Range("K16").Value = 1.392301088640352
$ws.Range("M16").Value = 0.4537164689676985
$ws.Range("N16").Value = 1.324975521471671
$ws.Range("O16").Value = 3.280840379743154
$ws.Range("C17").Value = 0.02975173689195287
$ws.Range("D17").Value = 0.1820286242274136
$ws.Range("E17").Value = 0.1469574221222096
$ws.Range("F17").Value = 1.373110443304711
$ws.Range("G17").Value = 0.002445218278220633
$ws.Range("I17").Value = 0.7126753481889452
$ws.Range("J17").Value = 0.1577467654044966
$ws.Range("K17").Value = 1.343117353947548
$ws.Range("M17").Value = 0.4427942083028711
$ws.Range("N17").Value = 1.327512202706018
$ws.Range("O17").Value = 3.280469680610878
$ws.Range("C18").Value = 0.02917207898184415
$ws.Range("D18").Value = 0.1814172617193748
$ws.Range("E18").Value = 0.1468079679514389
$ws.Range("F18").Value = 1.372474299846829
$ws.Range("G18").Value = 0.002445853389009267
$ws.Range("I18").Value = 0.7119863229052115
$ws.Range("J18").Value = 0.157831732206084
$ws.Range("K18").Value = 1.31482667405669
$ws.Range("M18").Value = 0.4365213519626963
$ws.Range("N18").Value = 1.329019062168456
$ws.Range("O18").Value = 3.280455522804147
$ws.Range("C19").Value = 0.02897565264794366
$ws.Range("D19").Value = 0.1812112746749648
$ws.Range("E19").Value = 0.1467584280230163
$ws.Range("F19").Value = 1.372272222554159
$ws.Range("G19").Value = 0.002446069952437986
$ws.Range("I19").Value = 0.7117606700929997
$ws.Range("J19").Value = 0.1578618336409221
$ws.Range("K19").Value = 1.30524773240893
$ws.Range("M19").Value = 0.43439908638085
$ws.Range("N19").Value = 1.329537477097631
$ws.Range("O19").Value = 3.280484904819303
$ws.Range("C20").Value = 0.02985894043392534
$ws.Range("D20").Value = 0.1821422528514063
$ws.Range("E20").Value = 0.1469855873489827
$ws.Range("F20").Value = 1.373234501789398
$ws.Range("G20").Value = 0.00244510145801046
$ws.Range("I20").Value = 0.7128065005464634
$ws.Range("J20").Value = 0.1577316732138101
$ws.Range("K20").Value = 1.348353214070414
$ws.Range("M20").Value = 0.4439559378360372
$ws.Range("N20").Value = 1.327237219584163
$ws.Range("O20").Value = 3.280488537283674
$ws.Range("C21").Value = 0.03281515265534551
$ws.Range("D21").Value = 0.1853401574913534
$ws.Range("E21").Value = 0.147822641901616
$ws.Range("F21").Value = 1.377393788807865
$ws.Range("G21").Value = 0.002441953330313873
$ws.Range("I21").Value = 0.7168521365773017
$ws.Range("J21").Value = 0.1573871034351058
$ws.Range("K21").Value = 1.493159894495477
$ws.Range("M21").Value = 0.4761763601716922
$ws.Range("N21").Value = 1.320079829445078
$ws.Range("O21").Value = 3.282880010935884
$ws.Range("C22").Value = 0.03473851324596922
$ws.Range("D22").Value = 0.1874819509652212
$ws.Range("E22").Value = 0.1484246203054518
$ws.Range("F22").Value = 1.380802704356157
$ws.Range("G22").Value = 0.002439974439937726
$ws.Range("I22").Value = 0.7198929674532906
$ws.Range("J22").Value = 0.1572310504781598
$ws.Range("K22").Value = 1.587774143217189
$ws.Range("M22").Value = 0.4973147947639234
$ws.Range("N22").Value = 1.31582713934668
$ws.Range("O22").Value = 3.286218453205947
$ws.Range("C23").Value = 0.03371278239931996
$ws.Range("D23").Value = 0.1863341204709457
$ws.Range("E23").Value = 0.1480983155717404
$ws.Range("F23").Value = 1.378920094540106
$ws.Range("G23").Value = 0.002441023440355829
$ws.Range("I23").Value = 0.7182336868149051
$ws.Range("J23").Value = 0.1573080113326526
$ws.Range("K23").Value = 1.537279630096521
$ws.Range("M23").Value = 0.4860255738197665
$ws.Range("N23").Value = 1.318058024038265
$ws.Range("O23").Value = 3.284274068199835
$ws.Range("C24").Value = 0.02981047751642052
$ws.Range("D24").Value = 0.1820908639405303
$ws.Range("E24").Value = 0.1469728348069488
$ws.Range("F24").Value = 1.373178174530693
$ws.Range("G24").Value = 0.002445154243959968
$ws.Range("I24").Value = 0.7127470690098008
$ws.Range("J24").Value = 0.15773847210491
$ws.Range("K24").Value = 1.345986126105402
$ws.Range("M24").Value = 0.4434306996482533
$ws.Range("N24").Value = 1.327361388405713
$ws.Range("O24").Value = 3.280479392452492
$ws.Range("C25").Value = 0.02556459974735503
$ws.Range("D25").Value = 0.177738398333986
$ws.Range("E25").Value = 0.145995308586123
$ws.Range("F25").Value = 1.369947475099025
$ws.Range("G25").Value = 0.002449949518975775
$ws.Range("I25").Value = 0.7085295518969374
$ws.Range("J25").Value = 0.1584995553702448
$ws.Range("K25").Value = 1.139594267160874
$ws.Range("M25").Value = 0.3978436784317836
$ws.Range("N25").Value = 1.339225658423814
$ws.Range("O25").Value = 3.28399474362871
